$wb = $excel.ActiveWorkbook

# --- Overview sheet: widen zh-cn / de-de status columns ---
# NB: Range.ColumnWidth is specified in characters and Excel rounds the
# stored OOXML <col width> to the nearest 1/6th of a character (pixel
# snapping), so we back the "characters" value off by the fixed 5/6
# padding that Excel always re-adds, landing on the desired stored width.
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Columns.Item(5).ColumnWidth = 29.144371396019366
$ovw.Columns.Item(6).ColumnWidth = 29.144371396019366

# The status text shown for both locales changes globally (shared string)
$ovw.Range("E2").Value = "Handed back: in sync with en-US"
$ovw.Range("F2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Columns.Item(3).ColumnWidth = 29.144371396019366
$zh.Columns.Item(9).ColumnWidth = 39.166666666666664
$zh.Columns.Item(10).ColumnWidth = 39.166666666666664

$zh.Range("C2").Value = "Handed back: in sync with en-US"

$zh.Range("I2").Value = "fa8c90cb-2def-4255-aea6-5cc38092d975.md"
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5e3bcaec526bfe454b7a6228610780b0ec4e9dbc/e2e/fa8c90cb-2def-4255-aea6-5cc38092d975.md", "", "", "fa8c90cb-2def-4255-aea6-5cc38092d975.md")

$zh.Range("J2").Value = "fa8c90cb-2def-4255-aea6-5cc38092d975.638e5c72b9ee0eaa35111d7d509c9183bbb2bc4a.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-17 03:01:15"

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Columns.Item(3).ColumnWidth = 29.144371396019366
$de.Columns.Item(9).ColumnWidth = 39.166666666666664
$de.Columns.Item(10).ColumnWidth = 39.166666666666664

$de.Range("C2").Value = "Handed back: in sync with en-US"

$de.Range("I2").Value = "fa8c90cb-2def-4255-aea6-5cc38092d975.md"
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5e3bcaec526bfe454b7a6228610780b0ec4e9dbc/e2e/fa8c90cb-2def-4255-aea6-5cc38092d975.md", "", "", "fa8c90cb-2def-4255-aea6-5cc38092d975.md")

$de.Range("J2").Value = "fa8c90cb-2def-4255-aea6-5cc38092d975.638e5c72b9ee0eaa35111d7d509c9183bbb2bc4a.de-de.xlf"
$de.Range("K2").Value = "2016-08-17 03:01:22"
